$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two incomplete rows (old rows 3 and 4: "361443-1"/"충남0908" and
# "4299-1"/"경기0831") so the remaining data rows (old rows 5-10) shift up to
# become rows 3-8.
$ws.Rows("3:4").Delete()

# Update the saved selection/active cell to match the post-merge cursor position.
$ws.Range("K19").Select()
